# Update the "settings" column (H) for every DNS record of type CNAME
# to reflect the new default export settings: {'flatten_cname': False}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $type = $ws.Cells.Item($r, 3).Value2
    if ($type -eq "CNAME") {
        $settings = $ws.Cells.Item($r, 8).Value2
        if ($settings -eq "{}") {
            $ws.Cells.Item($r, 8).Value = "{'flatten_cname': False}"
        }
    }
}
